# Weekly refresh of the "Puerro" (leek) price sheet:
# A brand-new observation is inserted at row 58 (pushing every existing
# row 58..145 down by one, to 59..146), and the previously-last row's
# data now lives at the newly created row 146. This mirrors the source
# diff, where every row's Fecha/Volumen/Precio.../Origen value equals
# what used to sit one row above it, and a single new record appears at
# the top of the block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 58:145 down to 59:146, leaving a blank row 58 behind
# (formats/styles of the row below are inherited automatically, same as
# Excel's native Insert behaviour).
$ws.Rows("58:58").Insert()

# Populate the new row 58 with this week's record.
$ws.Range("A58").Value = 10
$ws.Range("B58").Value = "Vega Modelo de Temuco"
$ws.Range("C58").Value = "La Araucanía"
$ws.Range("D58").Value = 44495
$ws.Range("E58").Value = 9
$ws.Range("F58").Value = 100112005
$ws.Range("G58").Value = "Puerro"
$ws.Range("H58").Value = "Azul de Maquehue"
$ws.Range("I58").Value = "Primera"
$ws.Range("J58").Value = 20
$ws.Range("K58").Value = 7000
$ws.Range("L58").Value = 7000
$ws.Range("M58").Value = 7000
$ws.Range("N58").Value = "`$/docena de paquetes"
$ws.Range("O58").Value = "Región de La Araucanía"
$ws.Range("P58").Value = 583
$ws.Range("Q58").Value = 12
$ws.Range("R58").Value = "Hortaliza"
